$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column widths (narrower columns)
$ws.Range("A1:A1").ColumnWidth = 26.833333333333332
$ws.Range("B1:B1").ColumnWidth = 29.5
$ws.Range("C1:C1").ColumnWidth = 28.0
$ws.Range("D1:D1").ColumnWidth = 34.333333333333336
$ws.Range("E1:E1").ColumnWidth = 28.166666666666668
$ws.Range("F1:F1").ColumnWidth = 26.833333333333332
$ws.Range("G1:G1").ColumnWidth = 33.166666666666664
$ws.Range("H1:H1").ColumnWidth = 28.333333333333332
$ws.Range("I1:I1").ColumnWidth = 27.166666666666668

# Update existing data rows 2-10 with new sweep results
$ws.Range("B2").Value = 34.904393155926961
$ws.Range("C2").Value = 16.089606651325976
$ws.Range("D2").Value = 0.46096222270501991
$ws.Range("E2").Value = 30.633069395226826
$ws.Range("F2").Value = 15.269849988791119
$ws.Range("G2").Value = 0.49847600290326871
$ws.Range("H2").Value = 339
$ws.Range("I2").Value = 301.5
$ws.Range("B3").Value = 34.90687719640686
$ws.Range("C3").Value = 16.092641777042623
$ws.Range("D3").Value = 0.46101636896637432
$ws.Range("E3").Value = 30.631262875553546
$ws.Range("F3").Value = 15.266611017638436
$ws.Range("G3").Value = 0.49839966049269685
$ws.Range("H3").Value = 339
$ws.Range("I3").Value = 301.5
$ws.Range("B4").Value = 34.902654088335453
$ws.Range("C4").Value = 16.090373407239497
$ws.Range("D4").Value = 0.46100715912653006
$ws.Range("E4").Value = 30.631330987920762
$ws.Range("F4").Value = 15.268912082004395
$ws.Range("G4").Value = 0.49847367350852551
$ws.Range("H4").Value = 339
$ws.Range("I4").Value = 301.5
$ws.Range("B5").Value = 34.900335217325917
$ws.Range("C5").Value = 16.083344130774595
$ws.Range("D5").Value = 0.46083637966864521
$ws.Range("E5").Value = 30.631006153298273
$ws.Range("F5").Value = 15.265590679661964
$ws.Range("G5").Value = 0.49837052701640366
$ws.Range("H5").Value = 338.5
$ws.Range("I5").Value = 301.5
$ws.Range("B6").Value = 34.902059062116741
$ws.Range("C6").Value = 16.084051007736623
$ws.Range("D6").Value = 0.46083387169539552
$ws.Range("E6").Value = 30.614430259734988
$ws.Range("F6").Value = 15.254588001013238
$ws.Range("G6").Value = 0.49828096984305231
$ws.Range("H6").Value = 338.5
$ws.Range("I6").Value = 301.5
$ws.Range("B7").Value = 34.913142317786942
$ws.Range("C7").Value = 16.085918286271752
$ws.Range("D7").Value = 0.46074106248742264
$ws.Range("E7").Value = 30.915778157333293
$ws.Range("F7").Value = 15.498930275205753
$ws.Range("G7").Value = 0.50132751620645755
$ws.Range("H7").Value = 338
$ws.Range("I7").Value = 301.5
$ws.Range("B8").Value = 34.954009431737539
$ws.Range("C8").Value = 16.131814095082728
$ws.Range("D8").Value = 0.46151541289095621
$ws.Range("E8").Value = 30.816875756889921
$ws.Range("F8").Value = 15.407115015040635
$ws.Range("G8").Value = 0.49995707340955775
$ws.Range("H8").Value = 337
$ws.Range("I8").Value = 301
$ws.Range("B9").Value = 34.754306317984621
$ws.Range("C9").Value = 15.974767991636696
$ws.Range("D9").Value = 0.45964859276647657
$ws.Range("E9").Value = 30.884770773729109
$ws.Range("F9").Value = 15.418543610312238
$ws.Range("G9").Value = 0.49922804100678003
$ws.Range("H9").Value = 334.5
$ws.Range("I9").Value = 299
$ws.Range("B10").Value = 34.997256415769407
$ws.Range("C10").Value = 16.07071326634804
$ws.Range("D10").Value = 0.45919923194627166
$ws.Range("E10").Value = 31.292030720013706
$ws.Range("F10").Value = 15.654514132322792
$ws.Range("G10").Value = 0.50027159542287236
$ws.Range("H10").Value = 326
$ws.Range("I10").Value = 294.5

# Append new row 11 (Hose thermal conductivity = 1000 data point)
$ws.Range("A11").Value = 1000
$ws.Range("B11").Value = 37.56518829049142
$ws.Range("C11").Value = 16.244581311393258
$ws.Range("D11").Value = 0.43243710601884888
$ws.Range("E11").Value = 32.310473124006961
$ws.Range("F11").Value = 15.92682582135299
$ws.Range("G11").Value = 0.49293075221232896
$ws.Range("H11").Value = 310
$ws.Range("I11").Value = 279
